$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue "D2" "306.40"
Set-TextValue "E2" "5.61%"
Set-TextValue "D3" "32.23"
Set-TextValue "E3" "9.30%"
Set-TextValue "D4" "5.340"
Set-TextValue "E4" "4.07%"
Set-TextValue "E5" "11.46%"
Set-TextValue "D6" "7.752"
Set-TextValue "E6" "5.32%"
Set-TextValue "D7" "3.701"
Set-TextValue "E7" "8.64%"
Set-TextValue "D8" "1.583"
Set-TextValue "E8" "16.76%"
Set-TextValue "D9" "0.9195"
Set-TextValue "E9" "0.39%"
Set-TextValue "D10" "0.01648"
Set-TextValue "E10" "2,455.13%"
Set-TextValue "D11" "0.1674"
Set-TextValue "D12" "0.07643"
Set-TextValue "E12" "14.40%"
Set-TextValue "D13" "0.07898"
Set-TextValue "E13" "2.26%"
Set-TextValue "D14" "0.03084"
Set-TextValue "E14" "4.82%"
Set-TextValue "D15" "0.09867"
Set-TextValue "E15" "9.67%"
Set-TextValue "D16" "0.001525"
Set-TextValue "E16" "-4.08%"
Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04555"
Set-TextValue "E17" "0.99%"
Set-TextValue "B18" "TigerCash"
Set-TextValue "C18" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.006381"
Set-TextValue "E18" "1.98%"
Set-TextValue "B19" "LEO"
Set-TextValue "C19" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D19" "3.473"
Set-TextValue "E19" "0.65%"
Set-TextValue "B20" "BTSEToken"
Set-TextValue "C20" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D20" "2.242"
Set-TextValue "E20" "0.93%"
Set-TextValue "B21" "BitpandaEcosystemToken"
Set-TextValue "C21" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D21" "0.3268"
Set-TextValue "E21" "1.68%"
Set-TextValue "B22" "ProBitToken"
Set-TextValue "C22" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D22" "0.1319"
Set-TextValue "E22" "0.71%"
Set-TextValue "B23" "MCDex"
Set-TextValue "C23" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D23" "4.214"
Set-TextValue "E23" "3.86%"
Set-TextValue "B24" "ZBToken"
Set-TextValue "C24" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D24" "0.1628"
Set-TextValue "E24" "4.98%"
Set-TextValue "D25" "0.001214"
Set-TextValue "E25" "1.80%"
Set-TextValue "D26" "0.004527"
Set-TextValue "E26" "9.54%"
Set-TextValue "D27" "0.0001169"
Set-TextValue "E27" "-6.43%"
Set-TextValue "D28" "0.0001742"
Set-TextValue "E28" "7.71%"
Set-TextValue "D40" "0.04512"
Set-TextValue "E40" "6.54%"
Set-TextValue "D41" "0.007432"
Set-TextValue "E41" "10.41%"
Set-TextValue "D42" "0.1368"
Set-TextValue "E42" "10.19%"
Set-TextValue "D43" "0.002258"
Set-TextValue "E43" "14.11%"
Set-TextValue "D44" "0.01374"
Set-TextValue "E44" "6.86%"
Set-TextValue "D45" "0.00006201"
Set-TextValue "E45" "10.21%"
Set-TextValue "D46" "1.892"
Set-TextValue "E46" "-4.11%"
Set-TextValue "D47" "0.01300"
Set-TextValue "E47" "-0.43%"
